# Updated main GSC export data.
#
# The GSC (Google Search Console) Coverage export rolled its date window
# forward by one day: the oldest date (2025-10-08) drops out of the
# "Chart" sheet's data table, and every subsequent row shifts up by one
# row (the table keeps the same trailing date, 2026-01-03, but is now
# one row shorter).
#
# The other sheets (Critical issues / Non-critical issues / Metadata)
# are unaffected content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 1 is the header (Date / Not indexed / Indexed / Impressions).
# Row 2 holds the oldest date (2025-10-08); remove it so every later
# row shifts up by one, matching the refreshed export window.
$ws.Rows.Item(2).Delete()
